$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '28.440.94'
$ws.Range('E2').Value = '  +0.43%  '

$ws.Range('D3').Value = '1.579.63'
$ws.Range('E3').Value = '  +0.13%  '

$ws.Range('E4').Value = '  -0.02%  '

$ws.Range('D5').Value = '''213.22'
$ws.Range('E5').Value = '  +0.68%  '

$ws.Range('E6').Value = '  +0.00%  '

$ws.Range('E7').Value = '  +0.00%  '

$ws.Range('D8').Value = '''44.69'
$ws.Range('E8').Value = '  -4.16%  '

$ws.Range('D9').Value = '''23.86'
$ws.Range('E9').Value = '  -0.14%  '

$ws.Range('E10').Value = '  -0.61%  '

$ws.Range('E11').Value = '  -1.09%  '

$ws.Range('E12').Value = '  +1.57%  '

$ws.Range('D13').Value = '1.805.58'
$ws.Range('E13').Value = '  +0.14%  '

$ws.Range('D14').Value = '1.580.76'
$ws.Range('E14').Value = '  +0.17%  '

$ws.Range('D15').Value = '''3.69'
$ws.Range('E15').Value = '  -0.94%  '

$ws.Range('B16').Value = 'WrappedBTC'
$ws.Range('C16').Value = 'https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc'
$ws.Range('D16').Value = '28.442.62'
$ws.Range('E16').Value = '  +0.33%  '

$ws.Range('B17').Value = 'Polygon'
$ws.Range('C17').Value = 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic'
$ws.Range('D17').Value = '''0.517'
$ws.Range('E17').Value = '  -1.71%  '

$ws.Range('D18').Value = '''61.93'
$ws.Range('E18').Value = '  -1.19%  '

$ws.Range('D19').Value = '''230.77'
$ws.Range('E19').Value = '  +1.20%  '

$ws.Range('E20').Value = '  +0.55%  '

$ws.Range('D21').Value = '0.0₃0687'

$ws.Range('E22').Value = '  +0.02%  '

$ws.Range('E23').Value = '  -0.16%  '

$ws.Range('D24').Value = '''9.10'
$ws.Range('E24').Value = '  -1.61%  '

$ws.Range('D25').Value = '''2.03'
$ws.Range('E25').Value = '  +0.89%  '

$ws.Range('D26').Value = '''151.71'
$ws.Range('E26').Value = '  +0.31%  '

$ws.Range('D27').Value = '''15.03'
$ws.Range('E27').Value = '  -0.36%  '

$ws.Range('D28').Value = '''6.41'
$ws.Range('E28').Value = '  -1.60%  '

$ws.Range('E29').Value = '  -0.86%  '

$ws.Range('E30').Value = '  -0.04%  '

$ws.Range('D31').Value = '''0.0482'
$ws.Range('E31').Value = '  +3.39%  '

$ws.Range('E32').Value = '  -1.52%  '

$ws.Range('E33').Value = '  -0.75%  '

$ws.Range('D34').Value = '''3.07'
$ws.Range('E34').Value = '  -1.77%  '

$ws.Range('D35').Value = '1.398.93'
$ws.Range('E35').Value = '  +0.78%  '

$ws.Range('D36').Value = '''1.09'
$ws.Range('E36').Value = '  +7.43%  '

$ws.Range('E37').Value = '  -3.91%  '

$ws.Range('D38').Value = '''2.36'
$ws.Range('E38').Value = '  +0.02%  '

$ws.Range('D39').Value = '''2.64'
$ws.Range('E39').Value = '  +1.90%  '

$ws.Range('E40').Value = '  -0.93%  '

$ws.Range('D41').Value = '''0.524'
$ws.Range('E41').Value = '  -2.86%  '

$ws.Range('E42').Value = '  +0.00%  '

$ws.Range('E43').Value = '  +1.27%  '

$ws.Range('D44').Value = '''0.787'
$ws.Range('E44').Value = '  -1.89%  '

$ws.Range('E45').Value = '  -3.06%  '

$ws.Range('E46').Value = '  -2.75%  '

$ws.Range('E47').Value = '  -5.42%  '

$ws.Range('D48').Value = '''62.67'
$ws.Range('E48').Value = '  -0.03%  '

$ws.Range('D49').Value = '1.717.73'
$ws.Range('E49').Value = '  +0.18%  '

$ws.Range('D50').Value = '''86.08'
$ws.Range('E50').Value = '  +0.02%  '

$ws.Range('D51').Value = '0.0₆0101'
$ws.Range('E51').Value = '  +2.20%  '
